# Updated cryptos list on Sat Apr 29 04:40:24 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free approach: for numeric-looking Price (column D) values we prefix
# with an apostrophe so Excel keeps them as literal text (matching the
# original inline-string cell type) instead of auto-converting to a number,
# then reset the cell style back to "Normal" so no stray NumberFormat/
# quotePrefix styling is left behind on the cell.

function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "29.524.59"
$ws.Range("E2").Value = "  +0.11%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.911.81"
$ws.Range("E3").Value = "  -0.15%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.55%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "326.17"
$ws.Range("E5").Value = "  -0.57%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  +0.52%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "0.4850"
$ws.Range("E7").Value = "  +1.26%  "

# Row 8 - Cardano
Set-TextValue $ws.Range("D8") "0.4077"
$ws.Range("E8").Value = "  -0.57%  "

# Row 9 - Dogecoin
Set-TextValue $ws.Range("D9") "0.08143"
$ws.Range("E9").Value = "  +1.43%  "

# Row 10 - Polygon
Set-TextValue $ws.Range("D10") "1.013"
$ws.Range("E10").Value = "  +0.15%  "

# Row 11 - Solana
Set-TextValue $ws.Range("D11") "23.49"
$ws.Range("E11").Value = "  +4.74%  "

# Row 12 - WrappedEther
Set-TextValue $ws.Range("D12") "1.899.34"
$ws.Range("E12").Value = "  -0.85%  "

# Row 13 - Polkadot
Set-TextValue $ws.Range("D13") "6.032"
$ws.Range("E13").Value = "  +1.23%  "

# Row 14 - Chainlink
Set-TextValue $ws.Range("D14") "7.107"
$ws.Range("E14").Value = "  -0.87%  "

# Row 15 - Litecoin
Set-TextValue $ws.Range("D15") "90.40"
$ws.Range("E15").Value = "  +1.12%  "

# Row 16 - was BinanceUSD, now TRON
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D16") "0.06790"
$ws.Range("E16").Value = "  +2.92%  "

# Row 17 - was TRON, now BinanceUSD
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D17") "1.007"
$ws.Range("E17").Value = "  +0.72%  "

# Row 18 - ShibaInu
Set-TextValue $ws.Range("D18") "0.00001041"
$ws.Range("E18").Value = "  +0.86%  "

# Row 19 - Avalanche
Set-TextValue $ws.Range("D19") "17.72"
$ws.Range("E19").Value = "  -0.35%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  +0.43%  "

# Row 21 - WrappedBTC
Set-TextValue $ws.Range("D21") "29.533.32"
$ws.Range("E21").Value = "  +0.08%  "

# Row 22 - Uniswap
Set-TextValue $ws.Range("D22") "5.611"
$ws.Range("E22").Value = "  +1.05%  "

# Row 23 - Cosmos
Set-TextValue $ws.Range("D23") "11.81"
$ws.Range("E23").Value = "  +2.31%  "

# Row 24 - Toncoin
Set-TextValue $ws.Range("D24") "2.167"
$ws.Range("E24").Value = "  -1.71%  "

# Row 25 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D25") "2.132.18"
$ws.Range("E25").Value = "  -0.61%  "

# Row 26 - Monero
Set-TextValue $ws.Range("D26") "154.86"
$ws.Range("E26").Value = "  +0.93%  "

# Row 27 - EthereumClassic
Set-TextValue $ws.Range("D27") "20.04"
$ws.Range("E27").Value = "  +1.14%  "

# Row 28 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D28") "6.279"
$ws.Range("E28").Value = "  +8.76%  "

# Row 29 - LidoDAOToken
Set-TextValue $ws.Range("D29") "2.107"
$ws.Range("E29").Value = "  -1.47%  "

# Row 30 - BitcoinCash
Set-TextValue $ws.Range("D30") "119.85"
$ws.Range("E30").Value = "  +2.04%  "

# Row 31 - ImmutableX
$ws.Range("E31").Value = "  -3.05%  "

# Row 32 - Stellar
Set-TextValue $ws.Range("D32") "0.09568"
$ws.Range("E32").Value = "  -0.08%  "

# Row 33 - Filecoin
Set-TextValue $ws.Range("D33") "5.540"
$ws.Range("E33").Value = "  +2.69%  "

# Row 34 - ARBITRUM
Set-TextValue $ws.Range("D34") "1.397"
$ws.Range("E34").Value = "  -1.97%  "

# Row 35 - HuobiToken
Set-TextValue $ws.Range("D35") "3.553"
$ws.Range("E35").Value = "  -0.53%  "

# Row 36 - VeChain
Set-TextValue $ws.Range("D36") "0.02267"
$ws.Range("E36").Value = "  +0.53%  "

# Row 37 - Hedera
Set-TextValue $ws.Range("D37") "0.06117"
$ws.Range("E37").Value = "  +0.26%  "

# Row 38 - TrustWalletToken
Set-TextValue $ws.Range("D38") "1.173"
$ws.Range("E38").Value = "  -0.13%  "

# Row 39 - TheSandbox
Set-TextValue $ws.Range("D39") "0.5948"
$ws.Range("E39").Value = "  +0.97%  "

# Row 40 - was FraxShare, now Aptos
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D40") "10.73"
$ws.Range("E40").Value = "  +5.68%  "

# Row 41 - was Aptos, now FraxShare
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D41") "7.921"
$ws.Range("E41").Value = "  -5.29%  "

# Row 42 - Algorand
Set-TextValue $ws.Range("D42") "0.1855"
$ws.Range("E42").Value = "  +0.78%  "

# Row 43 - RenderToken
Set-TextValue $ws.Range("D43") "2.446"
$ws.Range("E43").Value = "  +1.33%  "

# Row 44 - WEMIXToken
Set-TextValue $ws.Range("D44") "1.287"
$ws.Range("E44").Value = "  -0.96%  "

# Row 45 - Cronos
Set-TextValue $ws.Range("D45") "0.07723"
$ws.Range("E45").Value = "  -0.95%  "

# Row 46 - EnergySwap
Set-TextValue $ws.Range("D46") "12.39"
$ws.Range("E46").Value = "  +1.48%  "

# Row 47 - Decentraland
Set-TextValue $ws.Range("D47") "0.5575"
$ws.Range("E47").Value = "  +0.50%  "

# Row 48 - NEARProtocol
Set-TextValue $ws.Range("D48") "1.955"
$ws.Range("E48").Value = "  +1.21%  "

# Row 49 - Quant
Set-TextValue $ws.Range("D49") "114.85"
$ws.Range("E49").Value = "  +1.15%  "

# Row 50 - Aave
Set-TextValue $ws.Range("D50") "72.69"
$ws.Range("E50").Value = "  +1.37%  "

# Row 51 - EOS
Set-TextValue $ws.Range("D51") "1.053"
$ws.Range("E51").Value = "  +1.97%  "
